# Auto-generated edit script: apply numeric 'want-to-go' (F) / 'min price' (G)
# count corrections across the four sheets, plus one new row on 本地生活.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 416
$ws.Cells.Item(5, 6).Value = 57
$ws.Cells.Item(7, 6).Value = 39
$ws.Cells.Item(8, 6).Value = 1079
$ws.Cells.Item(10, 6).Value = 380
$ws.Cells.Item(13, 6).Value = 321
$ws.Cells.Item(14, 6).Value = 367
$ws.Cells.Item(15, 6).Value = 48
$ws.Cells.Item(18, 6).Value = 569
$ws.Cells.Item(19, 6).Value = 1471
$ws.Cells.Item(20, 6).Value = 5743
$ws.Cells.Item(21, 6).Value = 94
$ws.Cells.Item(22, 6).Value = 1617
$ws.Cells.Item(23, 6).Value = 382
$ws.Cells.Item(24, 6).Value = 64
$ws.Cells.Item(26, 6).Value = 5340
$ws.Cells.Item(27, 6).Value = 5340
$ws.Cells.Item(28, 6).Value = 133
$ws.Cells.Item(30, 6).Value = 1549
$ws.Cells.Item(33, 6).Value = 62
$ws.Cells.Item(34, 6).Value = 1050
$ws.Cells.Item(35, 6).Value = 669
$ws.Cells.Item(36, 6).Value = 111
$ws.Cells.Item(37, 6).Value = 1
$ws.Cells.Item(39, 6).Value = 3814
$ws.Cells.Item(17, 6).Value = 13
$ws.Cells.Item(17, 7).Value = 199

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 45
$ws.Cells.Item(8, 6).Value = 200
$ws.Cells.Item(5, 6).Value = 162
$ws.Cells.Item(5, 7).Value = 220

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 9419
$ws.Cells.Item(4, 6).Value = 2158

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 9419
$ws.Cells.Item(4, 6).Value = 2158
$ws.Cells.Item(6, 6).Value = 416
$ws.Cells.Item(8, 6).Value = 57
$ws.Cells.Item(10, 6).Value = 39
$ws.Cells.Item(11, 6).Value = 1079
$ws.Cells.Item(12, 6).Value = 380
$ws.Cells.Item(14, 6).Value = 321
$ws.Cells.Item(15, 6).Value = 367
$ws.Cells.Item(16, 6).Value = 48
$ws.Cells.Item(21, 6).Value = 1471
$ws.Cells.Item(22, 6).Value = 5743
$ws.Cells.Item(23, 6).Value = 94
$ws.Cells.Item(24, 6).Value = 1617
$ws.Cells.Item(27, 6).Value = 382
$ws.Cells.Item(30, 6).Value = 5340
$ws.Cells.Item(31, 6).Value = 5340
$ws.Cells.Item(32, 6).Value = 133
$ws.Cells.Item(34, 6).Value = 1549
$ws.Cells.Item(37, 6).Value = 1050
$ws.Cells.Item(38, 6).Value = 669
$ws.Cells.Item(39, 6).Value = 111
$ws.Cells.Item(47, 6).Value = 3814

# 本地生活: append new row 5 (新活动：剑网3×HAPPY ZOO)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 1).Copy()
$ws.Cells.Item(5, 1).PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).NumberFormat = "@"
$ws.Cells.Item(5, 2).Value = "2024-08-27"
$ws.Cells.Item(5, 3).Value = "杭州·剑网3×HAPPY ZOO 剑网3十五周年主题咖啡厅"
$ws.Cells.Item(5, 4).Value = "延安路292号（地铁1号线龙翔桥站D出口） 工联CC"
$ws.Cells.Item(5, 5).Value = "2024.08.27 00:00-10.07 23:59"
$ws.Cells.Item(5, 6).Value = 87
$ws.Cells.Item(5, 7).Value = 10
$ws.Cells.Item(5, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90962"
$ws.Cells.Item(5, 9).Value = "//i0.hdslb.com/bfs/openplatform/202408/tfJu8BDJ1724122581005.png"

